# Auto-generated edit script: updates numeric price/profit figures
# across multiple sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2421.05
$ws.Range("I32").Value = 3079.625
$ws.Range("J32").Value = 1982
$ws.Range("K32").Value = 3079.625
$ws.Range("L32").Value = 1982
$ws.Range("M32").Value = -2753.625
$ws.Range("N32").Value = -2634

$ws.Range("H98").Value = 37547.176
$ws.Range("J98").Value = 16324
$ws.Range("L98").Value = 16324
$ws.Range("N98").Value = -19320

$ws.Range("H113").Value = 7554.4614
$ws.Range("I113").Value = 9865.571
$ws.Range("J113").Value = 4858.1665
$ws.Range("K113").Value = 9865.571
$ws.Range("L113").Value = 4858.1665
$ws.Range("M113").Value = -6611.571
$ws.Range("N113").Value = -11366.1665

$ws.Range("H122").Value = 37547.176
$ws.Range("J122").Value = 16324
$ws.Range("L122").Value = 48972
$ws.Range("N122").Value = -53872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2680.4285
$ws.Range("I32").Value = 2645.5532
$ws.Range("K32").Value = 2645.5532
$ws.Range("M32").Value = -2358.5532

$ws.Range("H45").Value = 13482.875
$ws.Range("I45").Value = 20289.25
$ws.Range("K45").Value = 20289.25
$ws.Range("M45").Value = -19912.25

$ws.Range("H61").Value = 6621.4614
$ws.Range("I61").Value = 7229.778
$ws.Range("K61").Value = 7229.778
$ws.Range("M61").Value = -7017.778

$ws.Range("H132").Value = 2749.721
$ws.Range("I132").Value = 2500.2942
$ws.Range("K132").Value = 7500.882599999999
$ws.Range("M132").Value = -4970.882599999999

$ws.Range("H136").Value = 6621.4614
$ws.Range("I136").Value = 7229.778
$ws.Range("K136").Value = 21689.334
$ws.Range("M136").Value = -19139.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3055.8462
$ws.Range("I20").Value = 2147.889
$ws.Range("K20").Value = 2147.889
$ws.Range("M20").Value = -1900.889

$ws.Range("H86").Value = 6543.2666
$ws.Range("I86").Value = 11564.714
$ws.Range("J86").Value = 2149.5
$ws.Range("K86").Value = 11564.714
$ws.Range("L86").Value = 2149.5
$ws.Range("M86").Value = -10441.714
$ws.Range("N86").Value = -4395.5

$ws.Range("H89").Value = 6543.2666
$ws.Range("I89").Value = 11564.714
$ws.Range("J89").Value = 2149.5
$ws.Range("K89").Value = 57823.57
$ws.Range("L89").Value = 10747.5
$ws.Range("M89").Value = -52207.57
$ws.Range("N89").Value = -21979.5

$ws.Range("H99").Value = 12778.728
$ws.Range("I99").Value = 19284.428
$ws.Range("J99").Value = 1393.75
$ws.Range("K99").Value = 19284.428
$ws.Range("L99").Value = 1393.75
$ws.Range("M99").Value = -17786.428
$ws.Range("N99").Value = -4389.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 8092.846
$ws.Range("I10").Value = 10460.7
$ws.Range("K10").Value = 10460.7
$ws.Range("M10").Value = -10321.7

$ws.Range("H22").Value = 694.3333
$ws.Range("I22").Value = 665.4286
$ws.Range("K22").Value = 665.4286
$ws.Range("M22").Value = -315.4286

$ws.Range("H31").Value = 3492.125
$ws.Range("I31").Value = 2089.3157
$ws.Range("K31").Value = 2089.3157
$ws.Range("M31").Value = -1794.3157

$ws.Range("H34").Value = 3492.125
$ws.Range("I34").Value = 2089.3157
$ws.Range("K34").Value = 2089.3157
$ws.Range("M34").Value = -1887.3157

$ws.Range("H132").Value = 14526.955
$ws.Range("I132").Value = 6049.324
$ws.Range("K132").Value = 18147.972
$ws.Range("M132").Value = -15617.972

$ws.Range("H134").Value = 3132962.8
$ws.Range("I134").Value = 3915050.5
$ws.Range("J134").Value = 4611.75
$ws.Range("K134").Value = 11745151.5
$ws.Range("L134").Value = 13835.25
$ws.Range("M134").Value = -11742616.5
$ws.Range("N134").Value = -18905.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 32473234
$ws.Range("I4").Value = 31765482
$ws.Range("J4").Value = 44505000
$ws.Range("K4").Value = 95296446
$ws.Range("L4").Value = 133515000
$ws.Range("M4").Value = -95296334
$ws.Range("N4").Value = -133515224

$ws.Range("H38").Value = 1205.6364
$ws.Range("J38").Value = 2209.875
$ws.Range("L38").Value = 6629.625
$ws.Range("N38").Value = -7323.625

$ws.Range("H131").Value = 52632800
$ws.Range("J131").Value = 1771.1
$ws.Range("L131").Value = 5313.299999999999
$ws.Range("N131").Value = -15393.3

$ws.Range("H141").Value = 3293.3333
$ws.Range("I141").Value = 3052
$ws.Range("K141").Value = 9156
$ws.Range("M141").Value = -3976

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3558.4075
$ws.Range("I70").Value = 3173.4211
$ws.Range("J70").Value = 4472.75
$ws.Range("K70").Value = 3173.4211
$ws.Range("L70").Value = 4472.75
$ws.Range("M70").Value = -2903.4211
$ws.Range("N70").Value = -5012.75

$ws.Range("H73").Value = 3558.4075
$ws.Range("I73").Value = 3173.4211
$ws.Range("J73").Value = 4472.75
$ws.Range("K73").Value = 3173.4211
$ws.Range("L73").Value = 4472.75
$ws.Range("M73").Value = -2237.4211
$ws.Range("N73").Value = -6344.75

$ws.Range("H80").Value = 3571.2856
$ws.Range("I80").Value = 3667.3333
$ws.Range("J80").Value = 3499.25
$ws.Range("K80").Value = 3667.3333
$ws.Range("L80").Value = 3499.25
$ws.Range("M80").Value = -2669.3333
$ws.Range("N80").Value = -5495.25

$ws.Range("H83").Value = 3571.2856
$ws.Range("I83").Value = 3667.3333
$ws.Range("J83").Value = 3499.25
$ws.Range("K83").Value = 18336.6665
$ws.Range("L83").Value = 17496.25
$ws.Range("M83").Value = -13344.6665
$ws.Range("N83").Value = -27480.25

$ws.Range("H97").Value = 7560.636
$ws.Range("I97").Value = 9253.875
$ws.Range("K97").Value = 9253.875
$ws.Range("M97").Value = -8757.875

$ws.Range("H113").Value = 2814.6155
$ws.Range("I113").Value = 2298.5715
$ws.Range("K113").Value = 2298.5715
$ws.Range("M113").Value = -128.5715

$ws.Range("H119").Value = 60000
$ws.Range("J119").Value = 60000
$ws.Range("L119").Value = 60000
$ws.Range("N119").Value = -69676

$ws.Range("H120").Value = 64500
$ws.Range("J120").Value = 64500
$ws.Range("L120").Value = 64500
$ws.Range("N120").Value = -74176

$ws.Range("H126").Value = 16103.637
$ws.Range("I126").Value = 20488.334
$ws.Range("K126").Value = 61465.00199999999
$ws.Range("M126").Value = -58995.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1102607.2
$ws.Range("I132").Value = 1573579
$ws.Range("J132").Value = 3672.8667
$ws.Range("K132").Value = 4720737
$ws.Range("L132").Value = 11018.6001
$ws.Range("M132").Value = -4718207
$ws.Range("N132").Value = -16078.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 99675.05499999999
$ws.Range("I62").Value = 149896.95
$ws.Range("J62").Value = 3416.4167
$ws.Range("K62").Value = 149896.95
$ws.Range("L62").Value = 3416.4167
$ws.Range("M62").Value = -149272.95
$ws.Range("N62").Value = -4664.4167

$ws.Range("H65").Value = 99675.05499999999
$ws.Range("I65").Value = 149896.95
$ws.Range("J65").Value = 3416.4167
$ws.Range("K65").Value = 749484.75
$ws.Range("L65").Value = 17082.0835
$ws.Range("M65").Value = -746364.75
$ws.Range("N65").Value = -23322.0835

$ws.Range("H100").Value = 20007.58
$ws.Range("I100").Value = 12711.958
$ws.Range("J100").Value = 45021.145
$ws.Range("K100").Value = 25423.916
$ws.Range("L100").Value = 90042.28999999999
$ws.Range("M100").Value = -24882.916
$ws.Range("N100").Value = -91124.28999999999

$ws.Range("H132").Value = 8290.575000000001
$ws.Range("I132").Value = 13458.483
$ws.Range("J132").Value = 4476.1665
$ws.Range("K132").Value = 40375.449
$ws.Range("L132").Value = 13428.4995
$ws.Range("M132").Value = -37845.449
$ws.Range("N132").Value = -18488.4995
